$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 0.25
$ws.Range("B4").Value = 0.21875
$ws.Range("B5").Value = 0.21875
$ws.Range("B6").Value = 0.140625
$ws.Range("B9").Value = 0.171875
$ws.Range("B10").Value = 0.171875
$ws.Range("B11").Value = 0.1875
$ws.Range("B12").Value = 0.15625
$ws.Range("B13").Value = 0.15625
$ws.Range("B14").Value = 0.15625
$ws.Range("B15").Value = 0.15625
$ws.Range("B16").Value = 0.171875
$ws.Range("B17").Value = 0.171875
$ws.Range("B18").Value = 0.15625
$ws.Range("B19").Value = 0.15625
$ws.Range("B20").Value = 0.15625
$ws.Range("B21").Value = 0.15625
$ws.Range("B22").Value = 0.171875
$ws.Range("B23").Value = 0.171875
$ws.Range("B24").Value = 0.171875
$ws.Range("B25").Value = 0.171875
$ws.Range("B26").Value = 0.171875
$ws.Range("B27").Value = 0.171875
$ws.Range("B28").Value = 0.171875
$ws.Range("B29").Value = 0.171875
$ws.Range("B30").Value = 0.171875
$ws.Range("B37").Value = 0.140625
$ws.Range("B38").Value = 0.140625
$ws.Range("B39").Value = 0.140625
$ws.Range("B40").Value = 0.140625
$ws.Range("B41").Value = 0.140625
$ws.Range("B42").Value = 0.140625
$ws.Range("B43").Value = 0.140625
$ws.Range("B44").Value = 0.140625
$ws.Range("B66").Value = 0.140625
$ws.Range("B67").Value = 0.140625
$ws.Range("B68").Value = 0.140625
$ws.Range("B69").Value = 0.140625
$ws.Range("B70").Value = 0.140625
$ws.Range("B71").Value = 0.140625
$ws.Range("B72").Value = 0.140625
$ws.Range("B73").Value = 0.140625
$ws.Range("B74").Value = 0.140625
$ws.Range("B75").Value = 0.140625
$ws.Range("B76").Value = 0.140625
$ws.Range("B103").Value = 0.09375
$ws.Range("B104").Value = 0.125
$ws.Range("B105").Value = 0.1875
$ws.Range("B106").Value = 0.140625
$ws.Range("B107").Value = 0.078125
$ws.Range("B108").Value = 0.125
$ws.Range("B109").Value = 0.0625
$ws.Range("B110").Value = 0.15625
$ws.Range("B113").Value = 0.171875
$ws.Range("B115").Value = 0.109375
$ws.Range("B116").Value = 0.078125
$ws.Range("B117").Value = 0.09375

$newAddr = "<__main__.DisplayOutputs object at 0x7f7810052ac0>"
for ($r = 102; $r -le 118; $r++) {
    $ws.Range("A$r").Value = $newAddr
}
